$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows (A1:J87 -> A1:J89); this also expands the
# autoFilter range to match.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 88 picks up the same "banded" look previously used on row 72
# (the most recent row sharing that formatting), row 89 (the new last
# row) picks up the formatting that row 87 (the old last row) had.
$ws.Range("A72:J72").Copy()
$ws.Range("A88:J88").PasteSpecial(-4122)
$ws.Range("A87:J87").Copy()
$ws.Range("A89:J89").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data for 2020-06-06 and 2020-06-07.
$ws.Cells.Item(88, 1).Value = 43988
$ws.Cells.Item(88, 2).Value = 83105
$ws.Cells.Item(88, 3).Value = 229
$ws.Cells.Item(88, 4).Value = 1485
$ws.Cells.Item(88, 5).Value = 1
$ws.Cells.Item(88, 6).Value = 5
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 1
$ws.Cells.Item(88, 9).Value = 109
$ws.Cells.Item(88, 10).Value = 0

$ws.Cells.Item(89, 1).Value = 43989
$ws.Cells.Item(89, 2).Value = 83316
$ws.Cells.Item(89, 3).Value = 211
$ws.Cells.Item(89, 4).Value = 1485
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 5
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 109
$ws.Cells.Item(89, 10).Value = 0

# Match the saved selection state left by the author (last row selected).
$ws.Range("A89:J89").Select()
